# Apply VERMONT_2021 data-cleaning fixes:
#  1. Rename header columns to snake_case machine-readable names.
#  2. Title-case the lowercase connector words ("de", "del", "y") inside
#     place names so they read "De" / "Del" / "Y".
#  3. Drop the trailing footnote/metadata rows (62-66) and shrink the
#     worksheet dimension back down to A1:D60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row renames ---------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Place-name capitalization fixes ---------------------------------
$ws.Range("B2").Value  = "Comitán De Domínguez"
$ws.Range("B7").Value  = "San Cristóbal De Las Casas"
$ws.Range("A11").Value = "Ciudad De México"
$ws.Range("A14").Value = "Estado De México"
$ws.Range("B16").Value = "Naucalpan De Juárez"
$ws.Range("B17").Value = "San Felipe Del Progreso"
$ws.Range("B24").Value = "Nopala De Villagrán"
$ws.Range("B26").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B27").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B30").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B56").Value = "Martínez De La Torre"

# --- 3. Remove trailing footer rows (62-66) -----------------------------
$ws.Rows("62:66").Delete()
